$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Results")

$ws.Range("C16").Value = "C:\Users\COCO\onnxruntime_training_cuda_python\orttraining\orttraining\python\orttraining_pybind_state.cc:621 onnxruntime::python::addObjectMethodsForTraining::<lambda_6dd399ad6691adab5d0e0423ed8ce22d>::operator () [ONNXRuntimeError] : 1 : FAIL : Type Error: Type parameter (T) of Optype (Sub) bound to different types (tensor(float) and tensor(double) in node (onnx::Pow::60_Grad/Sub_1).`n"
$ws.Range("C39").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::179): A typestr: T, has unsupported type: tensor(bool)"
$ws.Range("C51").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::269): A typestr: T, has unsupported type: tensor(bool)"
$ws.Range("C52").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::271): A typestr: T, has unsupported type: tensor(bool)"
$ws.Range("C64").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::317): A typestr: T, has unsupported type: tensor(bool)"
$ws.Range("C65").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::319): A typestr: T, has unsupported type: tensor(bool)"
$ws.Range("C69").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::337): A typestr: T, has unsupported type: tensor(bool)"
$ws.Range("C70").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::339): A typestr: T, has unsupported type: tensor(bool)"
$ws.Range("C93").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::458): A typestr: T, has unsupported type: tensor(bool)"
$ws.Range("C102").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::495): X typestr: T, has unsupported type: tensor(uint8)"
$ws.Range("C223").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::966): X typestr: T, has unsupported type: tensor(uint8)"
$ws.Range("C239").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::971): X typestr: T, has unsupported type: tensor(uint8)"
